$wb = $excel.ActiveWorkbook

# "Pediatric Influenza Vaccine " sheet (sheet3)
$wsPed = $wb.Worksheets.Item("Pediatric Influenza Vaccine ")
$wsPed.Range("B3").Value = "Fluzone Pediatric dose No Preservative"
$wsPed.Range("B6").Value = "Fluarix Preservative-Free"
$wsPed.Range("B9").Value = "FluMist No Preservative"
$wsPed.Range("B10").Value = "Afluria No Preservative"
$wsPed.Range("H10").Value = "Merck (CSL product)"

# "Adult Influenza Vaccine " sheet (sheet4)
$wsAdult = $wb.Worksheets.Item("Adult Influenza Vaccine ")
$wsAdult.Range("B5").Value = "Agriflu No Preservative"
$wsAdult.Range("B7").Value = "Fluvirin Preservative-free"
$wsAdult.Range("B8").Value = "Fluarix Preservative-free"
$wsAdult.Range("B10").Value = "Flumist No Preservative"
